$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve existing "General" number format for column C while forcing text storage
$ws.Range("C2:C52").NumberFormat = "@"

$ws.Range("A2").Value = "Meitar Law Offices"
$ws.Range("B2").Value = "21sec"
$ws.Range("C2").Value = "1"
$ws.Range("A3").Value = "JSA"
$ws.Range("B3").Value = "15sec"
$ws.Range("C3").Value = "1"
$ws.Range("A4").Value = "Howse Williams"
$ws.Range("B4").Value = "20sec"
$ws.Range("C4").Value = "1"
$ws.Range("A5").Value = "GÖRG"
$ws.Range("B5").Value = "25sec"
$ws.Range("C5").Value = "1"
$ws.Range("A6").Value = "Mishcon Karas"
$ws.Range("B6").Value = "13sec"
$ws.Range("C6").Value = "1"
$ws.Range("A7").Value = "Aron Tadmor Levy"
$ws.Range("B7").Value = "23sec"
$ws.Range("C7").Value = "1"
$ws.Range("A8").Value = "Gianni And Origoni"
$ws.Range("B8").Value = "19sec"
$ws.Range("C8").Value = "1"
$ws.Range("A9").Value = "Harneys"
$ws.Range("B9").Value = "2min 5sec"
$ws.Range("C9").Value = "1"
$ws.Range("A10").Value = "Longan Law"
$ws.Range("B10").Value = "12sec"
$ws.Range("C10").Value = "1"
$ws.Range("A11").Value = "Ashurst"
$ws.Range("B11").Value = "24sec"
$ws.Range("C11").Value = "1"
$ws.Range("A12").Value = "Collas Crill"
$ws.Range("B12").Value = "13sec"
$ws.Range("C12").Value = "1"
$ws.Range("A13").Value = "Meyer Köring"
$ws.Range("B13").Value = "10sec"
$ws.Range("C13").Value = "1"
$ws.Range("A14").Value = "Grandall"
$ws.Range("B14").Value = "2min 27sec"
$ws.Range("C14").Value = "1"
$ws.Range("A15").Value = "Szecskay"
$ws.Range("B15").Value = "32sec"
$ws.Range("C15").Value = "1"
$ws.Range("A16").Value = "Lex Caribbean"
$ws.Range("B16").Value = "11sec"
$ws.Range("C16").Value = "1"
$ws.Range("A17").Value = "Borenius"
$ws.Range("B17").Value = "22sec"
$ws.Range("C17").Value = "1"
$ws.Range("A18").Value = "Magnusson Law"
$ws.Range("B18").Value = "14sec"
$ws.Range("C18").Value = "1"
$ws.Range("A19").Value = "Fox And Mandal"
$ws.Range("B19").Value = "14sec"
$ws.Range("C19").Value = "0"
$ws.Range("A20").Value = "Dentons"
$ws.Range("B20").Value = "2min 46sec"
$ws.Range("C20").Value = "1"
$ws.Range("A21").Value = "Havel Partners"
$ws.Range("B21").Value = "12sec"
$ws.Range("C21").Value = "1"
$ws.Range("A22").Value = "Paul Hastings"
$ws.Range("B22").Value = "13sec"
$ws.Range("C22").Value = "1"
$ws.Range("A23").Value = "Huiye Law"
$ws.Range("B23").Value = "10sec"
$ws.Range("C23").Value = "1"
$ws.Range("A24").Value = "Appleby Global"
$ws.Range("B24").Value = "25sec"
$ws.Range("C24").Value = "1"
$ws.Range("A25").Value = "LEX Logmannsstofa"
$ws.Range("B25").Value = "7sec"
$ws.Range("C25").Value = "1"
$ws.Range("A26").Value = "Mourant"
$ws.Range("B26").Value = "11sec"
$ws.Range("C26").Value = "1"
$ws.Range("A27").Value = "TC Law Firm"
$ws.Range("B27").Value = "8sec"
$ws.Range("C27").Value = "1"
$ws.Range("A28").Value = "BNT"
$ws.Range("B28").Value = "12sec"
$ws.Range("C28").Value = "1"
$ws.Range("A29").Value = "Guantao Law Firm"
$ws.Range("B29").Value = "23sec"
$ws.Range("C29").Value = "1"
$ws.Range("A30").Value = "Matheson"
$ws.Range("B30").Value = "10sec"
$ws.Range("C30").Value = "1"
$ws.Range("A31").Value = "Gide Loyrette Nouel"
$ws.Range("B31").Value = "17sec"
$ws.Range("C31").Value = "1"
$ws.Range("A32").Value = "DahlLaw"
$ws.Range("B32").Value = "21sec"
$ws.Range("C32").Value = "1"
$ws.Range("A33").Value = "Spencer West"
$ws.Range("B33").Value = "17sec"
$ws.Range("C33").Value = "1"
$ws.Range("A34").Value = "Kromann Reumert"
$ws.Range("B34").Value = "15sec"
$ws.Range("C34").Value = "1"
$ws.Range("A35").Value = "Hankun Law"
$ws.Range("B35").Value = "7sec"
$ws.Range("C35").Value = "1"
$ws.Range("A36").Value = "Schoenherr"
$ws.Range("B36").Value = "51sec"
$ws.Range("C36").Value = "1"
$ws.Range("A37").Value = "Arnold And Porter"
$ws.Range("B37").Value = "8sec"
$ws.Range("C37").Value = "1"
$ws.Range("A38").Value = "Gornitzky And Co"
$ws.Range("B38").Value = "10sec"
$ws.Range("C38").Value = "1"
$ws.Range("A39").Value = "Carey Olsen"
$ws.Range("B39").Value = "25sec"
$ws.Range("C39").Value = "1"
$ws.Range("A40").Value = "Legance"
$ws.Range("B40").Value = "16sec"
$ws.Range("C40").Value = "1"
$ws.Range("A41").Value = "Anand And Anand"
$ws.Range("B41").Value = "11sec"
$ws.Range("C41").Value = "1"
$ws.Range("A42").Value = "ALMT Legal"
$ws.Range("B42").Value = "15sec"
$ws.Range("C42").Value = "1"
$ws.Range("A43").Value = "Barnea And Co"
$ws.Range("B43").Value = "13sec"
$ws.Range("C43").Value = "1"
$ws.Range("A44").Value = "Byrne Wallace"
$ws.Range("B44").Value = "10sec"
$ws.Range("C44").Value = "1"
$ws.Range("A45").Value = "Beauchamps"
$ws.Range("B45").Value = "7sec"
$ws.Range("C45").Value = "0"
$ws.Range("A46").Value = "Dittmar And Indrenius"
$ws.Range("B46").Value = "10sec"
$ws.Range("C46").Value = "1"
$ws.Range("A47").Value = "Keystone Law"
$ws.Range("B47").Value = "16sec"
$ws.Range("C47").Value = "1"
$ws.Range("A48").Value = "Pedersoli"
$ws.Range("B48").Value = "8sec"
$ws.Range("C48").Value = "1"
$ws.Range("A49").Value = "Fangda Partners"
$ws.Range("B49").Value = "24sec"
$ws.Range("C49").Value = "1"
$ws.Range("A50").Value = "Fischer"
$ws.Range("B50").Value = "24sec"
$ws.Range("C50").Value = "1"
$ws.Range("A51").Value = "Deacons"
$ws.Range("B51").Value = "16sec"
$ws.Range("C51").Value = "1"
$ws.Range("A52").Value = "Myers Fletcher And Gordon"
$ws.Range("B52").Value = "17sec"
$ws.Range("C52").Value = "1"

# Restore original "General" number format for column C
$ws.Range("C2:C52").NumberFormat = "General"
